$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PL")

$ws.Range("B2").Value = 0.77710092075167148
$ws.Range("B3").Value = 0.80966819962611325
$ws.Range("B4").Value = 0.40685676867032922
$ws.Range("B5").Value = 0.48603505601304164
$ws.Range("B6").Value = 1.0366294171387864
$ws.Range("B7").Value = 1.0108885888661292
